# Chapter 3.1.22: "Ingen dokumentlyter er registrert." -> "Ingen dokumentflyter er registrert."
# The canonical edit splits the sentence's single run into three runs
# ("Ingen dokument" / "f" / "lyter er registrert.") -- reproduce that by
# inserting the missing "f" in place, then forcing a run boundary around
# the inserted character via a transient bookmark (added then immediately
# removed), which leaves no formatting residue but still breaks the run.

$d = $word.ActiveDocument

$needle = "Ingen dokumentlyter er registrert."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith($needle)) {
        $target = $p.Range
        break
    }
}

if ($target -ne $null) {
    $insertAt = $target.Start + ("Ingen dokument").Length
    $ins = $d.Range($insertAt, $insertAt)
    $ins.InsertAfter("f")

    $split = $d.Range($insertAt, $insertAt + 1)
    $d.Bookmarks.Add("tmp_run_split", $split)
    $d.Bookmarks("tmp_run_split").Delete()
}

Write-Output $target.Text
